# Apply the changes described by the diff to SORMAS_User_Rights.xlsx:
#  - Add three new user rights (rows) to the "User Rights" sheet:
#       LINE_LISTING_CONFIGURE, LINE_LISTING_CONFIGURE_NATION, AGGREGATE_REPORT_VIEW
#  - Bump the version string on the "About" sheet from 1.30.0-SNAPSHOT to 1.31.0-SNAPSHOT

$wb = $excel.ActiveWorkbook

$wsUserRights = $wb.Worksheets.Item("User Rights")
$wsAbout = $wb.Worksheets.Item("About")

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

$cols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

# New user rights to append to the "User Rights" sheet, starting at row 96.
# Each entry: Name, then Yes/No flags for columns C..V (20 role columns).
$newRights = @(
    @{
        Name = "LINE_LISTING_CONFIGURE"
        Flags = @("Yes","Yes","Yes","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No")
    },
    @{
        Name = "LINE_LISTING_CONFIGURE_NATION"
        Flags = @("Yes","Yes","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No","No")
    },
    @{
        Name = "AGGREGATE_REPORT_VIEW"
        Flags = @("Yes","Yes","Yes","Yes","Yes","No","Yes","No","Yes","No","No","No","No","Yes","Yes","Yes","Yes","Yes","Yes","Yes")
    }
)

# Update the version number on the "About" sheet first, so the renamed
# (pre-existing) strings keep sorting ahead of the brand new ones below.
$wsAbout.Range("A2").Value = "1.31.0-SNAPSHOT"

$startRow = 96

for ($i = 0; $i -lt $newRights.Count; $i++) {
    $row = $startRow + $i
    $entry = $newRights[$i]

    # Column A: bold user right name (same formatting as the existing rows, e.g. row 2).
    $wsUserRights.Range("A2").Copy()
    $wsUserRights.Range("A$row").PasteSpecial($xlPasteFormats)
    $wsUserRights.Range("A$row").Value = $entry.Name

    # Column B: plain user right name (no special style, as in existing rows).
    $wsUserRights.Range("B$row").Value = $entry.Name

    for ($c = 0; $c -lt $cols.Count; $c++) {
        $colLetter = $cols[$c]
        $flag = $entry.Flags[$c]
        $targetCell = $wsUserRights.Range($colLetter + $row)

        if ($flag -eq "Yes") {
            # Reuse the existing "Yes" (green) formatting, e.g. from C2.
            $wsUserRights.Range("C2").Copy()
        } else {
            # Reuse the existing "No" (red) formatting, e.g. from O2.
            $wsUserRights.Range("O2").Copy()
        }
        $targetCell.PasteSpecial($xlPasteFormats)
        $targetCell.Value = $flag
    }
}

$excel.CutCopyMode = $false

$wb.Save()
